$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Format as text first so the date-like string "10/06/2025" is stored
# literally (matching the existing rows) rather than being auto-converted
# into a date serial number by Excel's smart input parsing.
$ws.Range("A50").NumberFormat = "@"
$ws.Range("A50").Value = "10/06/2025"
$ws.Range("A50").ClearFormats()

$ws.Range("B50").Value = 15060.38
